$d = $word.ActiveDocument

# 1. Merge the "Annex7" + "1" runs into a single "Annex71" run.
#    Word's Find/Replace naturally coalesces the matched text into one run,
#    carrying over the formatting of the first matched run (Consolas font).
$d.Content.Find.Execute("Annex71", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Annex71", 2)

# 2. Remove the trailing empty paragraph at the end of the document by
#    deleting the paragraph mark that separates it from the preceding
#    paragraph (equivalent to pressing Delete at the end of the
#    second-to-last paragraph). This merges the two paragraphs, keeping
#    the formatting of the earlier (surviving) paragraph, and matches the
#    diff which drops the last, empty <w:p>.
$paras = $d.Paragraphs
$count = $paras.Count
$lastPara = $paras.Item($count)
$prevPara = $paras.Item($count - 1)
$mergeRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
$mergeRange.Delete()
